# Fix #6770 - Modify nutrient and packagingMaterial exports
#
# On the "Eco emballage" sheet, the hidden header row (row 2) stores the
# BIRT/export column formulas as text. Two of those formulas are
# simplified:
#   D2: excel|IF(ISBLANK(C2),IF(ISBLANK(B2),"",B2),C2)   ->  excel|IF(C2="",B2,C2)
#   G2: excel|IF(ISBLANK(F2),IF(ISBLANK(E2),"",E2), F2)  ->  excel|IF(F2="",E2,F2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eco emballage")

# --- D2 ---------------------------------------------------------------
# Replace the old ISBLANK()-based formula text with the simplified one.
$cellD2 = $ws.Range("D2")
$cellD2.Value = 'excel|IF(C2="",B2,C2)'

# Keep a bit of the former rich-text coloring on the tail of the string
# (matches how the cell looked after the manual edit in Excel/Calc).
$cellD2.Characters(13, 9).Font.Color = 0

# --- G2 -----------------------------------------------------------------
$ws.Range("G2").Value = 'excel|IF(F2="",E2,F2)'

# --- Selection ----------------------------------------------------------
# The active cell on the sheet ended up on I15 after the edit.
$ws.Range("I15").Select() | Out-Null
